{"js": "// Apply the README/docx stats fix-up for the Renaissance/JDK21/ShenandoahGC\n// finagle-http 4G benchmark table.\n//\n// The document is a single table where (almost) every row has exactly one\n// cell holding one numeric/text value. We:\n//   1. Rewrite the first four data rows (indices 0-3).\n//   2. Insert three new one-cell rows right after row 3.\n//   3. Rewrite three more rows (now shifted by the insertion).\n//   4. Delete three rows that are no longer needed.\n//   5. Rewrite the row that follows the deleted block.\n//   6. Collapse the three multi-value (tab-separated) summary rows near the\n//      end of the table down to a single short value each.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\n// --- 1. First four rows -----------------------------------------------\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\ntable.getCell(3, 0).value = \"1176\";\nawait context.sync();\n\n// --- 2. Insert three new rows right after row index 3 ------------------\ntable.rows.load(\"items\");\nawait context.sync();\n\ntable.rows.items[3].insertRows(\"After\", 3, [\n  [\"0.00001\"],\n  [\"0.00058\"],\n  [\"0.00012\"],\n]);\nawait context.sync();\n\n// --- 3. Rewrite the next three rows (now at indices 8, 9, 10) ----------\ntable.getCell(8, 0).value = \"0.00020\";\ntable.getCell(9, 0).value = \"0.00021\";\ntable.getCell(10, 0).value = \"0.00026\";\nawait context.sync();\n\n// --- 4. Delete the three now-unneeded rows (indices 11, 12, 13) --------\n// Delete from the highest index down so earlier indices stay valid.\ntable.rows.load(\"items\");\nawait context.sync();\ntable.rows.items[13].delete();\ntable.rows.items[12].delete();\ntable.rows.items[11].delete();\nawait context.sync();\n\n// --- 5. Rewrite the row that follows the deleted block (now index 11) --\ntable.getCell(11, 0).value = \"0.16178\";\nawait context.sync();\n\n// --- 6. Collapse the three trailing multi-value summary rows -----------\n// These are the last three rows of the table and are unaffected by the\n// insert/delete operations performed above (which all happened near the\n// start of the table), so we can address them from the end.\ntable.rows.load(\"items\");\nawait context.sync();\nconst rowCount = table.rows.items.length;\n\ntable.getCell(rowCount - 3, 0).value = \"99.96\";\ntable.getCell(rowCount - 2, 0).value = \"0.16\";\ntable.getCell(rowCount - 1, 0).value = \"381\";\nawait context.sync();\n", "ps1": "# Apply the README/docx stats fix-up for the Renaissance/JDK21/ShenandoahGC\n# finagle-http 4G benchmark table.\n#\n# The document is a single table where (almost) every row has exactly one\n# cell holding one numeric/text value. We:\n#   1. Rewrite the first four data rows (rows 1-4, 1-based).\n#   2. Insert three new one-cell rows right after row 4.\n#   3. Rewrite three more rows (now shifted by the insertion).\n#   4. Delete three rows that are no longer needed.\n#   5. Rewrite the row that follows the deleted block.\n#   6. Collapse the three multi-value (tab-separated) summary rows near the\n#      end of the table down to a single short value each.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- 1. First four rows -------------------------------------------------\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"1176\"\n\n# --- 2. Insert three new rows right after row 4 --------------------------\n# Rows.Add(beforeRow) inserts immediately before the given row, so we keep\n# targeting the row right after the previous insertion to preserve order.\n$newRow1 = $t.Rows.Add($t.Rows.Item(5))\n$newRow1.Cells.Item(1).Range.Text = \"0.00001\"\n$newRow2 = $t.Rows.Add($t.Rows.Item(6))\n$newRow2.Cells.Item(1).Range.Text = \"0.00058\"\n$newRow3 = $t.Rows.Add($t.Rows.Item(7))\n$newRow3.Cells.Item(1).Range.Text = \"0.00012\"\n\n# --- 3. Rewrite the next three rows (now at rows 9, 10, 11) -------------\n$t.Cell(9, 1).Range.Text = \"0.00020\"\n$t.Cell(10, 1).Range.Text = \"0.00021\"\n$t.Cell(11, 1).Range.Text = \"0.00026\"\n\n# --- 4. Delete the three now-unneeded rows (rows 12, 13, 14) ------------\n# Delete from the highest index down so earlier indices stay valid.\n$t.Rows.Item(14).Delete()\n$t.Rows.Item(13).Delete()\n$t.Rows.Item(12).Delete()\n\n# --- 5. Rewrite the row that follows the deleted block (now row 12) -----\n$t.Cell(12, 1).Range.Text = \"0.16178\"\n\n# --- 6. Collapse the three trailing multi-value summary rows ------------\n# These are the last three rows of the table and are unaffected by the\n# insert/delete operations performed above (which all happened near the\n# start of the table), so we can address them relative to the row count.\n$n = $t.Rows.Count\n$t.Cell($n - 2, 1).Range.Text = \"99.96\"\n$t.Cell($n - 1, 1).Range.Text = \"0.16\"\n$t.Cell($n, 1).Range.Text = \"381\"\n"}
